# ------------------------------------------------------------------
# Applies the "cost.xlsx" update:
#  1. Cost sheet gains a new pricing/discount lookup layout:
#       - "wac_per_unit" renamed to "price_per_unit" (col L, unchanged values)
#       - two new columns inserted before "infusion_cost": discount_lower,
#         discount_upper (uniform 0.2 / 0.3 sampled-uniform-distribution bounds)
#       - two new columns appended after "infusion_cost": loading_dose,
#         weight_based (0/1 flags per therapeutic agent)
#       - abatacept SC (row 8) now carries an infusion_cost of 164 and is
#         flagged loading_dose = 1; infliximab (row 4) is flagged
#         weight_based = 1
#  2. Lookup sheet: nbt's agent1 lookup corrected from "nbt" to "cdmards"
#  3. The Lookup sheet becomes the active/visible tab
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cost")

# --- Insert two blank columns ahead of the old "infusion_cost" column (M) ---
# Before: ... K=strength_unit, L=wac_per_unit, M=infusion_cost
# After:  ... K=strength_unit, L=price_per_unit, M=discount_lower,
#          N=discount_upper, O=infusion_cost, P=loading_dose, Q=weight_based
$ws.Range("M1:N1").EntireColumn.Insert()

# Rename the price column header
$ws.Range("L1").Value = "price_per_unit"

# New headers
$ws.Range("M1").Value = "discount_lower"
$ws.Range("N1").Value = "discount_upper"
$ws.Range("P1").Value = "loading_dose"
$ws.Range("Q1").Value = "weight_based"

# --- Discount bounds: uniform-distribution sampling bounds, same for every agent ---
$ws.Range("M2:M14").Value = 0.2
$ws.Range("N2:N14").Value = 0.3
$ws.Range("M2:N14").NumberFormat = "0.00"

# --- infusion_cost (moved to col O) keeps its old values, except abatacept SC ---
$ws.Range("O2").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("O4").Value = 164
$ws.Range("O5").Value = 0
$ws.Range("O6").Value = 0
$ws.Range("O7").Value = 164
$ws.Range("O8").Value = 164
$ws.Range("O9").Value = 0
$ws.Range("O10").Value = 164
$ws.Range("O11").Value = 0
$ws.Range("O12").Value = 0
$ws.Range("O13").Value = 0
$ws.Range("O14").Value = 0

# --- loading_dose flag (col P): only abatacept SC requires an IV loading dose ---
$ws.Range("P2:P14").Value = 0
$ws.Range("P8").Value = 1

# --- weight_based flag (col Q): only infliximab is dosed by body weight (mg/kg) ---
$ws.Range("Q2:Q14").Value = 0
$ws.Range("Q4").Value = 1

$ws.Range("P2:Q14").NumberFormat = "0"
$ws.Range("O2:O14").NumberFormat = "0"

$ws.Columns.Item("A:Q").AutoFit()

# --- Lookup sheet: fix the agent1 lookup for the "nbt" row ---
$wsLookup = $wb.Worksheets.Item("Lookup")
$wsLookup.Range("B15").Value = "cdmards"

# --- Make "Lookup" the active/visible sheet (matches workbookView activeTab=1) ---
$wsLookup.Activate()
$wsLookup.Select()
